# Typos from teaching the class.
#
# Colors several code-identifier runs (and the stray leading space that
# used to sit in front of the following word) red (FF0000) across a few
# slides, and tidies up the now-redundant leading space in the run that
# follows each recolored identifier.
#
# Helper: locate a substring inside a paragraph's text (PowerPoint's own
# 1-based Paragraph.Start plus .NET's 0-based IndexOf) and return the
# matching sub-TextRange, always carving it out of the *top level*
# TextRange for the shape (slicing via a previously-sliced TextRange
# does not use absolute document coordinates).
# NB: default parameter values aren't honored by this host, so every
# call below passes $startFrom explicitly (normally 0).
function Get-SubRange($fullRange, $para, $search, $startFrom) {
    $t = $para.Text
    $idx = $t.IndexOf($search, $startFrom)
    if ($idx -lt 0) {
        throw "substring not found: [$search]"
    }
    $globalStart = $para.Start + $idx
    return $fullRange.Characters($globalStart, $search.Length)
}

function Set-Red($range) {
    $range.Font.Color.RGB = 255
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 14 - "Adds a req.csrfToken() function to retrieve the token..."
# ---------------------------------------------------------------------
$s = $p.Slides.Item(14)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)

Set-Red (Get-SubRange $tr $para "req.csrfToken" 0)
Set-Red (Get-SubRange $tr $para "()" 0)
# Recolor the leading space of " function to " red; this splits it into
# its own run and leaves "function to " (no leading space) behind.
Set-Red (Get-SubRange $tr $para " function to " 0)

# ---------------------------------------------------------------------
# Slide 16 - "compression should be high in the middleware stack..."
# ---------------------------------------------------------------------
$s = $p.Slides.Item(16)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$para = $tr.Paragraphs(1, 1)

Set-Red (Get-SubRange $tr $para "compression" 0)
Set-Red (Get-SubRange $tr $para " should be" 0)
$resEnd = Get-SubRange $tr $para "res.end" 0
Set-Red $resEnd
$afterResEnd = ($resEnd.Start - $para.Start) + $resEnd.Length
Set-Red (Get-SubRange $tr $para "()" $afterResEnd)

# ---------------------------------------------------------------------
# Slide 4 - "...populate req.cookies with an object keyed..."
# ---------------------------------------------------------------------
$s = $p.Slides.Item(4)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$para = $tr.Paragraphs(1, 1)

Set-Red (Get-SubRange $tr $para "req.cookies" 0)

# ---------------------------------------------------------------------
# Slide 5 - "Populates req.body property" / "...inflates gzip and deflate encodings"
# ---------------------------------------------------------------------
$s = $p.Slides.Item(5)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

$para6 = $tr.Paragraphs(6, 1)
Set-Red (Get-SubRange $tr $para6 "req.body" 0)

$para7 = $tr.Paragraphs(7, 1)
Set-Red (Get-SubRange $tr $para7 "gzip" 0)
Set-Red (Get-SubRange $tr $para7 " and " 0)
Set-Red (Get-SubRange $tr $para7 "deflate" 0)
Set-Red (Get-SubRange $tr $para7 " encodings" 0)

# ---------------------------------------------------------------------
# Slide 6 - Body-Parser API: bodyParser.json / .raw / .text / .urlencoded
# ---------------------------------------------------------------------
$s = $p.Slides.Item(6)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

$para1 = $tr.Paragraphs(1, 1)
$para4 = $tr.Paragraphs(4, 1)
$para6 = $tr.Paragraphs(6, 1)
$para8 = $tr.Paragraphs(8, 1)

Set-Red (Get-SubRange $tr $para1 "bodyParser.json" 0)
Set-Red (Get-SubRange $tr $para4 "bodyParser.raw" 0)
Set-Red (Get-SubRange $tr $para6 "bodyParser.text" 0)
Set-Red (Get-SubRange $tr $para8 "bodyParser.urlencoded" 0)
